$d = $word.ActiveDocument

# --- Region 1 -----------------------------------------------------------
# Before: "Um einen Hamming Code zu erstellen muss zunächst die der Frame
#          "Initalize Hamming Code" ausgewählt werden. ..."
# After:  "Um einen Hamming-Code zu erstellen, muss zunächst der Frame
#          "Initialize Hamming Code" ausgewählt werden. ..."
#
# Four separate, minimally-scoped edits so each Find/Replace stays inside
# a single existing run and does not bridge the "Hamming" / quote runs
# (which are wrapped in their own <w:proofErr> spell-check markers).

# 1) " Code zu erstellen " -> "-Code zu erstellen, "
#    (turns "Hamming Code" into "Hamming-Code" and adds the comma after
#    "erstellen")
$r1 = $d.Content
$r1.Find.Execute(" Code zu erstellen ", $false, $false, $false, $false, `
    $false, $true, 1, $false, "-Code zu erstellen, ", 2) | Out-Null

# 2) "muss zunächst die der Frame " -> "muss zunächst der Frame "
#    (drops the stray "die ")
$r2 = $d.Content
$r2.Find.Execute("muss zunächst die der Frame ", $false, $false, $false, `
    $false, $false, $true, 1, $false, "muss zunächst der Frame ", 2) | Out-Null

# 3) "Initalize" -> "Initialize" (fix the typo)
$r3 = $d.Content
$r3.Find.Execute("Initalize", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Initialize", 2) | Out-Null
